$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'27.794.27"
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Formula = "'1.807.71"
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("D4").Formula = "'0.9990"
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").Formula = "'306.82"
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").Formula = "'0.9977"
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("D7").Formula = "'0.4986"
$ws.Range("E7").Value = '  -4.09%  '
$ws.Range("D8").Formula = "'0.3887"
$ws.Range("E8").Value = '  +2.33%  '
$ws.Range("D9").Formula = "'0.09477"
$ws.Range("E9").Value = '  +21.25%  '
$ws.Range("D10").Formula = "'1.097"
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").Formula = "'40.46"
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("D12").Formula = "'6.309"
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").Formula = "'0.9978"
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Formula = "'20.58"
$ws.Range("D15").Formula = "'1.798.29"
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Formula = "'7.207"
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").Formula = "'0.00001120"
$ws.Range("E17").Value = '  +3.57%  '
$ws.Range("D18").Formula = "'92.66"
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").Formula = "'0.06566"
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").Formula = "'0.9995"
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Formula = "'17.09"
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Formula = "'5.918"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Formula = "'27.861.91"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Formula = "'10.98"
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").Formula = "'2.221"
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Formula = "'156.61"
$ws.Range("E26").Value = '  -2.51%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Formula = "'20.59"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").Formula = "'2.412"
$ws.Range("E28").Value = '  +4.16%  '
$ws.Range("D29").Formula = "'2.004.40"
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("D30").Formula = "'127.16"
$ws.Range("E30").Value = '  +4.04%  '
$ws.Range("D31").Formula = "'0.1070"
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").Formula = "'1.055"
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").Formula = "'5.546"
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("D34").Formula = "'3.596"
$ws.Range("E34").Value = '  -2.13%  '
$ws.Range("D35").Formula = "'0.06793"
$ws.Range("E35").Value = '  -5.67%  '
$ws.Range("D36").Formula = "'8.848"
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").Formula = "'0.2134"
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Formula = "'11.37"
$ws.Range("E39").Value = '  -6.54%  '
$ws.Range("D40").Formula = "'4.928"
$ws.Range("E40").Value = '  -2.34%  '
$ws.Range("D41").Formula = "'0.6187"
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("D42").Formula = "'0.9981"
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Formula = "'1.142"
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").Formula = "'0.5865"
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Formula = "'1.273"
$ws.Range("E46").Value = '  -7.18%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Formula = "'3.662"
$ws.Range("E47").Value = '  -2.52%  '
$ws.Range("D48").Formula = "'123.36"
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("D49").Formula = "'1.940"
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("E50").Value = '  -3.98%  '
$ws.Range("D51").Formula = "'0.06727"
$ws.Range("E51").Value = '  +0.01%  '
